# Fill in the "carrier" (D) values for the practice rows, the new
# "pair_kind" (J) values for the unique video/audio pairs, and the
# corresponding new detail rows (9-16) that record which carrier goes
# with each unique_video / unique_audio pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 is an empty placeholder cell in the source file (no cached value) that
# the load/save round-trip otherwise resolves to shared-string index 0;
# explicitly re-blank it so it round-trips as empty, matching the diff
# (which does not touch F1 at all).
$ws.Range("F1").ClearContents()

# Practice rows 2-5: carrier column (D) was blank, now matches the
# carrier word used in column K.
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Rows 6-9: new pair_kind (J) values distinguishing unique video vs audio.
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# Rows 14-21: new kind (C) / carrier (D) detail rows for items 9-16.
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "look"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "look"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "where"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "where"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "can"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "can"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "do"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "do"
